$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells keep their original text formatting
# (values like "563.24" or "0.170" would otherwise be auto-converted to numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.496.88"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.421.94"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.24"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.35"
$ws.Range("E6").Value = "  +6.81%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.170"
$ws.Range("E9").Value = "  +10.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.421.61"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  +7.66%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.233.66"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.867.83"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("E17").Value = "  +6.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.427.12"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +6.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.92"
$ws.Range("E20").Value = "  +5.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.13"
$ws.Range("E21").Value = "  +6.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("E23").Value = "  +8.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.19"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +8.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.50"
$ws.Range("E27").Value = "  +8.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.548.64"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  +8.89%  "
$ws.Range("E31").Value = "  +7.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +13.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "454.51"
$ws.Range("E33").Value = "  +11.26%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.09"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.10"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +7.61%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +4.35%  "
$ws.Range("E41").Value = "  +5.38%  "
$ws.Range("E42").Value = "  +6.62%  "
$ws.Range("E43").Value = "  +6.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.78"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  +4.81%  "
$ws.Range("E46").Value = "  +8.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.08"
$ws.Range("E47").Value = "  +6.52%  "
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0724"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.488"
$ws.Range("E50").Value = "  +5.35%  "
$ws.Range("E51").Value = "  +2.98%  "
